$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - force text to avoid numeric auto-conversion,
# then clear formats so no extra style/numFmt is introduced.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.332.64"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.013.79"
$ws.Range("D3").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.77"
$ws.Range("D5").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.80"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.384"
$ws.Range("D9").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.33"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.311.13"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.66"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.797"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.22"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.003.13"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.334.48"
$ws.Range("D18").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0837"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "233.59"
$ws.Range("D21").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.60"
$ws.Range("D23").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.96"
$ws.Range("D26").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.61"
$ws.Range("D28").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.32"
$ws.Range("D30").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.60"
$ws.Range("D32").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.48"
$ws.Range("D34").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.81"
$ws.Range("D36").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.49"
$ws.Range("D39").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.19"
$ws.Range("D41").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.435.55"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.74"
$ws.Range("D45").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.55"
$ws.Range("D46").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.02"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.202.30"
$ws.Range("D50").ClearFormats()

# Volume(1h) (column E) updates
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E5").Value = "  +4.80%  "
$ws.Range("E6").Value = "  -1.95%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -7.01%  "
$ws.Range("E9").Value = "  -3.38%  "
$ws.Range("E10").Value = "  -5.00%  "
$ws.Range("E12").Value = "  -7.02%  "
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("E14").Value = "  -3.85%  "
$ws.Range("E15").Value = "  -7.87%  "
$ws.Range("E16").Value = "  -6.18%  "
$ws.Range("E17").Value = "  -1.65%  "
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("E20").Value = "  -3.95%  "
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("E22").Value = "  -3.17%  "
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("E27").Value = "  -5.93%  "
$ws.Range("E28").Value = "  -1.64%  "
$ws.Range("E29").Value = "  -5.97%  "
$ws.Range("E30").Value = "  -4.82%  "
$ws.Range("E31").Value = "  -2.03%  "
$ws.Range("E32").Value = "  -5.37%  "
$ws.Range("E33").Value = "  -5.10%  "
$ws.Range("E34").Value = "  -1.68%  "
$ws.Range("E35").Value = "  -5.49%  "
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("E38").Value = "  -8.42%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  +2.48%  "
$ws.Range("E41").Value = "  -0.63%  "
$ws.Range("E44").Value = "  +3.09%  "
$ws.Range("E45").Value = "  -8.61%  "
$ws.Range("E46").Value = "  -3.90%  "
$ws.Range("E47").Value = "  -3.87%  "
$ws.Range("E48").Value = "  +2.35%  "
$ws.Range("E49").Value = "  -6.88%  "
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("E51").Value = "  -10.61%  "

# Row 42/43: VeChain and Cronos swap places, with new Volume(1h) values
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0929"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.78%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0212"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.80%  "
